# Updates the Price (column D) values (and two Volume/text cells in column E)
# in the "cryptos" worksheet to reflect freshly scraped data, per the
# "Updated symbol list" GitHub Actions commit.
#
# The Price column is stored as text (e.g. "245.29"), even though the
# values look numeric, so each cell value is forced to remain text
# (matching the original inlineStr string cells) rather than letting Excel
# auto-convert the numeric-looking string into a real number. Re-applying
# the "Normal" cell style after the write keeps the cell's style/format
# untouched (Excel otherwise stamps a "Text" quote-prefixed number format
# on any cell whose text needs disambiguating from a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "245.42"
Set-TextValue "D3"  "24.12"
Set-TextValue "D4"  "5.278"
Set-TextValue "D5"  "0.05777"
Set-TextValue "D6"  "6.460"
Set-TextValue "D7"  "3.129"
Set-TextValue "D8"  "0.8163"
Set-TextValue "D9"  "0.8506"
Set-TextValue "D10" "0.1351"
Set-TextValue "D11" "0.06933"
Set-TextValue "D12" "0.03132"
Set-TextValue "D13" "0.02937"
Set-TextValue "D14" "0.09391"
Set-TextValue "D15" "3.753"
Set-TextValue "D17" "0.04663"
Set-TextValue "D18" "0.0005966"
Set-TextValue "D19" "0.006236"
Set-TextValue "D20" "0.001238"
Set-TextValue "D21" "0.004618"
Set-TextValue "D22" "0.00006895"
Set-TextValue "D23" "3.506"
Set-TextValue "D24" "2.138"
Set-TextValue "D25" "0.3194"
Set-TextValue "D26" "0.1320"
Set-TextValue "D27" "0.1358"
Set-TextValue "D28" "0.0002330"
Set-TextValue "D40" "0.03629"
Set-TextValue "D41" "0.006226"
Set-TextValue "E41" "40KickTokenKICKBestin24h"
Set-TextValue "D42" "0.1052"
Set-TextValue "D43" "0.003398"
Set-TextValue "D44" "0.008402"
Set-TextValue "E44" "43LocalTradersLCT"
Set-TextValue "D45" "0.00005270"
Set-TextValue "D47" "0.3697"
Set-TextValue "D48" "0.002274"
